$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")
$ws.Range("E2").Value = "Product iphoneNeetu904 Successfully Added`n✖︎"
